$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55 (G55=5517)
$ws.Range("H55").Value = 1245.625
$ws.Range("I55").Value = 143.33333
$ws.Range("J55").Value = 1500
$ws.Range("K55").Value = 143.33333
$ws.Range("L55").Value = 1500
$ws.Range("M55").Value = 70.66667000000001
$ws.Range("N55").Value = -1928

# Row 80 (G80=12605)
$ws.Range("H80").Value = 6355.1113
$ws.Range("I80").Value = 3258.4
$ws.Range("J80").Value = 10226
$ws.Range("K80").Value = 9775.200000000001
$ws.Range("L80").Value = 30678
$ws.Range("M80").Value = -8777.200000000001
$ws.Range("N80").Value = -32674

# Row 83 (G83=12605)
$ws.Range("H83").Value = 6355.1113
$ws.Range("I83").Value = 3258.4
$ws.Range("J83").Value = 10226
$ws.Range("K83").Value = 29325.6
$ws.Range("L83").Value = 92034
$ws.Range("M83").Value = -24333.6
$ws.Range("N83").Value = -102018

# Row 88 (G88=12608)
$ws.Range("H88").Value = 1828.5714
$ws.Range("J88").Value = 1828.5714
$ws.Range("L88").Value = 1828.5714
$ws.Range("N88").Value = -2640.5714

# Row 91 (G91=12608)
$ws.Range("H91").Value = 1828.5714
$ws.Range("J91").Value = 1828.5714
$ws.Range("L91").Value = 1828.5714
$ws.Range("N91").Value = -4636.5714

# Row 94 (G94=19905)
$ws.Range("H94").Value = 992.6667
$ws.Range("I94").Value = 992.6667
$ws.Range("K94").Value = 992.6667
$ws.Range("M94").Value = -541.6667

# Row 116 (G116=27778)
$ws.Range("H116").Value = 6166.3335
$ws.Range("I116").Value = 5999.5
$ws.Range("J116").Value = 6500
$ws.Range("K116").Value = 5999.5
$ws.Range("L116").Value = 6500
$ws.Range("M116").Value = -2557.5
$ws.Range("N116").Value = -13384

# Row 132 (G132=44049)
$ws.Range("H132").Value = 2319.5122
$ws.Range("I132").Value = 1534.2858
$ws.Range("K132").Value = 4602.857400000001
$ws.Range("M132").Value = -2072.857400000001

# Row 137 (G137=44013)
$ws.Range("H137").Value = 2750
$ws.Range("I137").Value = 2750
$ws.Range("K137").Value = 8250
$ws.Range("M137").Value = -5700

# Row 138 (G138=44169)
$ws.Range("H138").Value = 2816.3076
$ws.Range("J138").Value = 2959.5
$ws.Range("L138").Value = 8878.5
$ws.Range("N138").Value = -19158.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (G32=44147)
$ws.Range("H32").Value = 2175.4595
$ws.Range("I32").Value = 2175.4595
$ws.Range("K32").Value = 2175.4595
$ws.Range("M32").Value = -1888.4595

# Row 45 (G45=27714)
$ws.Range("H45").Value = 5276.125
$ws.Range("I45").Value = 5276.125
$ws.Range("K45").Value = 5276.125
$ws.Range("M45").Value = -4899.125

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (G22=5367)
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Row 31 (G31=44023)
$ws.Range("H31").Value = 2144.5715
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

# Row 34 (G34=44023)
$ws.Range("H34").Value = 2144.5715
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

# Row 86 (G86=12584)
$ws.Range("H86").Value = 5685.6924
$ws.Range("I86").Value = 3358.7144
$ws.Range("K86").Value = 3358.7144
$ws.Range("M86").Value = -2235.7144

# Row 89 (G89=12584)
$ws.Range("H89").Value = 5685.6924
$ws.Range("I89").Value = 3358.7144
$ws.Range("K89").Value = 16793.572
$ws.Range("M89").Value = -11177.572

# Row 132 (G132=44019)
$ws.Range("H132").Value = 2287.2942
$ws.Range("I132").Value = 2256
$ws.Range("K132").Value = 6768
$ws.Range("M132").Value = -4238

# Row 134 (G134=44020)
$ws.Range("H134").Value = 674.0968
$ws.Range("I134").Value = 691.9
$ws.Range("J134").Value = 140
$ws.Range("K134").Value = 2075.7
$ws.Range("L134").Value = 420
$ws.Range("M134").Value = 459.3000000000002
$ws.Range("N134").Value = -5490

$ws = $wb.Worksheets.Item("CUL")
# Row 68 (G68=12895)
$ws.Range("H68").Value = 1412.6666
$ws.Range("I68").Value = 499
$ws.Range("K68").Value = 1497
$ws.Range("M68").Value = -686

# Row 71 (G71=12895)
$ws.Range("H71").Value = 1412.6666
$ws.Range("I71").Value = 499
$ws.Range("K71").Value = 4491
$ws.Range("M71").Value = -435

# Row 131 (G131=36060)
$ws.Range("H131").Value = 590264.2
$ws.Range("J131").Value = 911693.8
$ws.Range("L131").Value = 2735081.4
$ws.Range("N131").Value = -2745161.4

# Row 137 (G137=44088)
$ws.Range("H137").Value = 1468.125
$ws.Range("I137").Value = 1311.25
$ws.Range("J137").Value = 1625
$ws.Range("K137").Value = 3933.75
$ws.Range("L137").Value = 4875
$ws.Range("M137").Value = 1166.25
$ws.Range("N137").Value = -15075

$ws = $wb.Worksheets.Item("GSM")
# Row 33 (G33=4450)
$ws.Range("H33").Value = 60000000
$ws.Range("J33").Value = 60000000
$ws.Range("L33").Value = 60000000
$ws.Range("N33").Value = -60000504

# Row 80 (G80=12521)
$ws.Range("H80").Value = 1874.5
$ws.Range("I80").Value = 1874.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1874.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -876.5
$ws.Range("N80").ClearContents()

# Row 83 (G83=12521)
$ws.Range("H83").Value = 1874.5
$ws.Range("I83").Value = 1874.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 9372.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -4380.5
$ws.Range("N83").ClearContents()

# Row 122 (G122=36182)
$ws.Range("H122").Value = 3603.3572
$ws.Range("I122").Value = 3522.5454
$ws.Range("K122").Value = 10567.6362
$ws.Range("M122").Value = -8117.636200000001

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (G16=5289)
$ws.Range("H16").Value = 793.36365
$ws.Range("I16").Value = 857.7
$ws.Range("J16").Value = 150
$ws.Range("K16").Value = 857.7
$ws.Range("L16").Value = 150
$ws.Range("M16").Value = -687.7
$ws.Range("N16").Value = -490

$ws = $wb.Worksheets.Item("WVR")
# Row 39 (G39=3106)
$ws.Range("H39").Value = 20000
$ws.Range("I39").Value = 20000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 20000
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -19587
$ws.Range("N39").ClearContents()

# Row 70 (G70=11979)
$ws.Range("H70").Value = 7900
$ws.Range("I70").Value = 7900
$ws.Range("K70").Value = 7900
$ws.Range("M70").Value = -7585

# Row 73 (G73=11979)
$ws.Range("H73").Value = 7900
$ws.Range("I73").Value = 7900
$ws.Range("K73").Value = 7900
$ws.Range("M73").Value = -6808

# Row 107 (G107=27746)
$ws.Range("H107").Value = 1076.6428
$ws.Range("I107").Value = 1149.75
$ws.Range("J107").Value = 979.1667
$ws.Range("K107").Value = 3449.25
$ws.Range("L107").Value = 2937.5001
$ws.Range("M107").Value = -1529.25
$ws.Range("N107").Value = -6777.5001

# Row 110 (G110=25825)
$ws.Range("H110").Value = 47200.5
$ws.Range("J110").Value = 47200.5
$ws.Range("L110").Value = 47200.5
$ws.Range("N110").Value = -55380.5

# Row 132 (G132=44029)
$ws.Range("H132").Value = 10961.7
$ws.Range("I132").Value = 10961.7
$ws.Range("K132").Value = 32885.10000000001
$ws.Range("M132").Value = -30355.10000000001
